# Update "想去人数" (number of people interested) figures on the
# "展览" (Exhibitions) sheet and the "全部类型" (All Types) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1100
$ws1.Range("F4").Value = 1756
$ws1.Range("F5").Value = 781
$ws1.Range("F6").Value = 139

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1100
$ws4.Range("F4").Value = 1756
$ws4.Range("F6").Value = 781
$ws4.Range("F7").Value = 139
